{"js": "// Office.js (Word JavaScript API) edit script.\n// Body of `async (context) => { ... }`.\n\n// 1) Bump the activation date: 01/01/2023 -> 01/01/2024.\nconst dateResults = context.document.body.search(\"Ativa\\u00e7\\u00e3o: 01/01/2023\", { matchCase: true });\ndateResults.load(\"text\");\nawait context.sync();\ndateResults.items[0].insertText(\"Ativa\\u00e7\\u00e3o: 01/01/2024\", Word.InsertLocation.replace);\n\n// 2) Drop the \"1341653 - Maria Jos\u00e9 Ramos Sandim\" line (and its line break)\n//    from the \"Docente(s) Respons\u00e1vel(eis)\" list.\nconst nameResults = context.document.body.search(\"1341653 - Maria Jos\u00e9 Ramos Sandim\\u000b\", { matchCase: true });\nnameResults.load(\"text\");\nawait context.sync();\nnameResults.items[0].insertText(\"\", Word.InsertLocation.replace);\nawait context.sync();\n\n// Re-fetch paragraphs once the list edit above has settled, then rewrite the\n// \"Programa resumido\" / \"Programa\" paragraphs (PT + italic EN) by index \u2014\n// this keeps existing run formatting (e.g. the <w:i/> italics) intact.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst texts = paragraphs.items.map((p) => p.text);\n\n// 3) \"Programa resumido\" section (Portuguese, then italic English).\nconst resumoPtIndex = texts.indexOf(\"\\u00d3ptica de raios; Ondas eletromagn\\u00e9ticas: fase e polariza\\u00e7\\u00e3o; Interfer\\u00eancia; Coer\\u00eancia; Difra\\u00e7\\u00e3o; \\u00d3ptica de Fourier; Intera\\u00e7\\u00e3o da luz com a mat\\u00e9ria; Guias de ondas met\\u00e1licos e diel\\u00e9tricos; \\u00d3ptica de cristais; \\u00d3ptica n\\u00e3o linear.\");\nparagraphs.items[resumoPtIndex].insertText(\"Descri\\u00e7\\u00e3o ondulat\\u00f3ria e qu\\u00e2ntica da luz. Propriedades da luz. Intera\\u00e7\\u00e3o da luz com a mat\\u00e9ria. Aplica\\u00e7\\u00f5es.\", Word.InsertLocation.replace);\n\nconst resumoEnIndex = texts.indexOf(\"Ray optics; Electromagnetic waves: phase and polarization; Interference; Coherence; Diffraction; Fourier optics; Interaction of light with matter; Metallic and dielectric waveguides; Crystal optics; Non-linear optics.\");\nparagraphs.items[resumoEnIndex].insertText(\"Presentation of the wave and quantum description of light, study of the properties of light, the interaction of light with matter and applications of physical optics.\", Word.InsertLocation.replace);\n\n// 4) \"Programa\" section (Portuguese, then italic English).\nconst programaPtIndex = texts.indexOf(\"\\u00d3ptica de raios. Introdu\\u00e7\\u00e3o. Propaga\\u00e7\\u00e3o de luz em meios homog\\u00eaneos. Propaga\\u00e7\\u00e3o de luz em meios n\\u00e3o homog\\u00eaneos. A lei de Snell generalizada. O princ\\u00edpio de Fermat. A equa\\u00e7\\u00e3o dos raios. A fun\\u00e7\\u00e3o eikonal. Analogia ente a mec\\u00e2nica cl\\u00e1ssica e a \\u00f3ptica geom\\u00e9trica. O potencial \\u00f3ptico.Ondas eletromagn\\u00e9ticas. Ondas harm\\u00f4nicas unidimensionais. Ondas planas e esf\\u00e9ricas. Ondas gaussianas. Propaga\\u00e7\\u00e3o do feixe gaussiano. Vetor de Poynting. Intensidade.A fase da onda eletromagn\\u00e9tica. Velocidades de fase e de grupo. Dispers\\u00e3o. Efeito Doppler. Aplica\\u00e7\\u00f5es astron\\u00f4micas. Alargamento de linhas espectrais. \\u00d3ptica relativ\\u00edstica. Modula\\u00e7\\u00e3o eletro\\u00f3ptica de frequ\\u00eancia. Automodula\\u00e7\\u00e3o de fase. Polariza\\u00e7\\u00e3o das ondas eletromagn\\u00e9ticas. Polariza\\u00e7\\u00e3o linear. Polariza\\u00e7\\u00e3o el\\u00edptica. Polariza\\u00e7\\u00e3o circular. Obten\\u00e7\\u00e3o de luz linearmente polarizada. Equa\\u00e7\\u00f5es de Fresnel. Polariza\\u00e7\\u00e3o por reflex\\u00e3o total interna. Matrizes de Jones. Atividade \\u00f3ptica. Efeito Faraday. Isoladores \\u00f3pticos. Efeito Pockels. Efeitos Kerr e Cotton-Mouton. Chaveamento eletro\\u00f3ptico.Interfer\\u00eancia. Princ\\u00edpio da superposi\\u00e7\\u00e3o. Interfer\\u00eancia por divis\\u00e3o da frente de onda. Interfer\\u00eancia por divis\\u00e3o de amplitudes. Interfer\\u00f4metro de Fabry-Perot. Analisador de espectro \\u00f3ptico. Teoria de pel\\u00edculas.Coer\\u00eancia. Introdu\\u00e7\\u00e3o. Coer\\u00eancia temporal. Resolu\\u00e7\\u00e3o espectral de um trem de ondas finito. Coer\\u00eancia espacial. Medidas de di\\u00e2metros de estrelas.Difra\\u00e7\\u00e3o. Princ\\u00edpio de Huygens. F\\u00f3rmula de Fresnel-Kirchhoff. Princ\\u00edpio de Babinet. Difra\\u00e7\\u00e3o de Fraunhofer. Difra\\u00e7\\u00e3o por uma abertura circular. Rede de difra\\u00e7\\u00e3o. Padr\\u00f5es de difra\\u00e7\\u00e3o de Fresnel. \\u00d3ptica de Fourier.  Microscopia por contraste de fase.  Holografia. Intera\\u00e7\\u00e3o da radia\\u00e7\\u00e3o com a mat\\u00e9ria.  Modelo do oscilador harm\\u00f4nico.  Dispers\\u00e3o crom\\u00e1tica do \\u00edndice de refra\\u00e7\\u00e3o. Absor\\u00e7\\u00e3o. Espalhamento Rayleigh. For\\u00e7a da radia\\u00e7\\u00e3o em \\u00e1tomo neutro.\\u00d3ptica n\\u00e3o linear. Susceptibilidade n\\u00e3o linear, processos param\\u00e9tricos e n\\u00e3o param\\u00e9tricos. Gera\\u00e7\\u00e3o de freq\\u00fc\\u00eancias. Casamento de fase.\");\nparagraphs.items[programaPtIndex].insertText(\"O que \\u00e9 luz? Reflex\\u00e3o. Refra\\u00e7\\u00e3o. Difra\\u00e7\\u00e3o. Polariza\\u00e7\\u00e3o. Forma\\u00e7\\u00e3o de imagens: Transformada de Fourier. Ondas eletromagn\\u00e9ticas. Equa\\u00e7\\u00f5es de Maxwell. Propaga\\u00e7\\u00e3o da luz em diferentes meios: v\\u00e1cuo, diel\\u00e9trico, condutor. Transporte de energia. Condi\\u00e7\\u00f5es de contorno entre diferentes meios: v\\u00e1cuo, diel\\u00e9trico, condutor. Propaga\\u00e7\\u00e3o da luz entre diferentes meios: incid\\u00eancia normal e obl\\u00edqua na interface entre meios. Coeficientes de Fresnel. Aplica\\u00e7\\u00f5es da \\u00d3ptica: holografia, laser, fibras \\u00f3pticas, materiais eletrocr\\u00f4micos, metamateriais.\", Word.InsertLocation.replace);\n\nconst programaEnIndex = texts.indexOf(\"Ray optics. Introduction. Propagation of light in homogeneous media. Propagation of light in non-homogeneous media. Generalized Snell's law. Fermat's principle. The equation of rays. The eikonal function. Analogy between classical mechanics and optics geometric The optical potential.Electromagnetic waves. One-dimensional harmonic waves. Flat and spherical waves. Gaussian waves. Gaussian beam propagation. Poynting vector. Intensity.The phase of the electromagnetic wave. Phase and group speeds. Dispersal. Doppler effect. Astronomical applications. Broadening of spectral lines. Relativistic optics. Electro-optical frequency modulation. Phase automodulation.Polarization of electromagnetic waves. Linear polarization. Elliptical Polarization. Circular polarization. Obtaining linearly polarized light. Fresnel equations. Polarization by total internal reflection. Jones matrices. Optical activity. Faraday effect. Optical isolators. Pockels Effect. Kerr and Cotton-Mouton effects. Electro-optical switching.Interference. Superposition principle. Interference by division of the wavefront. Amplitude division interference. Fabry-Perot interferometer. Optical spectrum analyzer. Film theory.Coherence. Introduction. Temporal coherence. Spectral resolution of a finite wave train. Spatial coherence. Star diameter measurements.Diffraction. Huygens Principle. Fresnel-Kirchhoff formula. Babinet's Principle. Fraunhofer Diffraction. Diffraction through a circular aperture. Diffraction grating. Fresnel diffraction patterns. Fourier optics. Phase contrast microscopy. Holography.Interaction of radiation with matter. Harmonic oscillator model. Chromatic dispersion of the refractive index. Absorption. Rayleigh scattering. Force of radiation on a neutral atom.Non-linear optics. Nonlinear susceptibility, parametric and nonparametric processes. Frequency generation. Phase marriage.\");\nparagraphs.items[programaEnIndex].insertText(\"What is light? Reflection. Refraction. Diffraction. Polarization. Image formation: Fourier transform. Electromagnetic waves. Maxwell's equations. Propagation of light in different media: vacuum, dielectric, conductor. Energy transport. Boundary conditions between different media: vacuum, dielectric, conductor. Light propagation between different media: normal and oblique incidence at the interface between media. Fresnel coefficients. Optics Applications: holography, laser, optical fibers, electrochromic materials, metamaterials.\", Word.InsertLocation.replace);\n\n// 5) Append a new reference to the \"Bibliografia\" paragraph.\nconst biblioIndex = texts.indexOf(\"HECHT, E.; ZAJAC, A. Optics; Reading, Addison-Wesley, 1974.ZILLIO, S. C. \\u00d3ptica Moderna - Fundamentos e Aplica\\u00e7\\u00f5es, 2005.\");\nparagraphs.items[biblioIndex].insertText(\"J. R. Reitz, F. J. Milford, R. W. Christy, Fundamentos da Teoria Eletromagn\\u00e9tica. Editora Campus. 1982.\", Word.InsertLocation.end);\n\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# $word.ActiveDocument is already open as $d.\n$d = $word.ActiveDocument\n\n# 1) Bump the activation date: 01/01/2023 -> 01/01/2024.\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"Ativa\u00e7\u00e3o: 01/01/2023\"\n$find.Replacement.Text = \"Ativa\u00e7\u00e3o: 01/01/2024\"\n$find.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $null, 2) | Out-Null\n\n# 2) Drop the \"1341653 - Maria Jos\u00e9 Ramos Sandim\" line (and its line break)\n#    from the \"Docente(s) Respons\u00e1vel(eis)\" list. \"^l\" matches the manual\n#    line break (<w:br/>) that follows the run.\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"1341653 - Maria Jos\u00e9 Ramos Sandim^l\"\n$find.Replacement.Text = \"\"\n$find.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $null, 2) | Out-Null\n\n# 3)-6) Rewrite the \"Programa resumido\" / \"Programa\" paragraphs (PT, then\n#    italic EN) by assigning Range.Text directly -- this preserves existing\n#    run formatting (e.g. the italic run) and, unlike Find.Execute's\n#    Replacement.Text, does not run the straight quotes through smart-quote\n#    autocorrection (relevant for \"Maxwell's equations\" below).\nfunction Set-ParagraphText($paragraphs, [string]$oldText, [string]$newText) {\n    for ($i = 1; $i -le $paragraphs.Count; $i++) {\n        $p = $paragraphs.Item($i)\n        if ($p.Range.Text.TrimEnd(\"`r\", \"`a\") -eq $oldText) {\n            $p.Range.Text = $newText\n            return $true\n        }\n    }\n    return $false\n}\n\n$paragraphs = $d.Paragraphs\nSet-ParagraphText $paragraphs \"\u00d3ptica de raios; Ondas eletromagn\u00e9ticas: fase e polariza\u00e7\u00e3o; Interfer\u00eancia; Coer\u00eancia; Difra\u00e7\u00e3o; \u00d3ptica de Fourier; Intera\u00e7\u00e3o da luz com a mat\u00e9ria; Guias de ondas met\u00e1licos e diel\u00e9tricos; \u00d3ptica de cristais; \u00d3ptica n\u00e3o linear.\" \"Descri\u00e7\u00e3o ondulat\u00f3ria e qu\u00e2ntica da luz. Propriedades da luz. Intera\u00e7\u00e3o da luz com a mat\u00e9ria. Aplica\u00e7\u00f5es.\" | Out-Null\nSet-ParagraphText $paragraphs \"Ray optics; Electromagnetic waves: phase and polarization; Interference; Coherence; Diffraction; Fourier optics; Interaction of light with matter; Metallic and dielectric waveguides; Crystal optics; Non-linear optics.\" \"Presentation of the wave and quantum description of light, study of the properties of light, the interaction of light with matter and applications of physical optics.\" | Out-Null\nSet-ParagraphText $paragraphs \"\u00d3ptica de raios. Introdu\u00e7\u00e3o. Propaga\u00e7\u00e3o de luz em meios homog\u00eaneos. Propaga\u00e7\u00e3o de luz em meios n\u00e3o homog\u00eaneos. A lei de Snell generalizada. O princ\u00edpio de Fermat. A equa\u00e7\u00e3o dos raios. A fun\u00e7\u00e3o eikonal. Analogia ente a mec\u00e2nica cl\u00e1ssica e a \u00f3ptica geom\u00e9trica. O potencial \u00f3ptico.Ondas eletromagn\u00e9ticas. Ondas harm\u00f4nicas unidimensionais. Ondas planas e esf\u00e9ricas. Ondas gaussianas. Propaga\u00e7\u00e3o do feixe gaussiano. Vetor de Poynting. Intensidade.A fase da onda eletromagn\u00e9tica. Velocidades de fase e de grupo. Dispers\u00e3o. Efeito Doppler. Aplica\u00e7\u00f5es astron\u00f4micas. Alargamento de linhas espectrais. \u00d3ptica relativ\u00edstica. Modula\u00e7\u00e3o eletro\u00f3ptica de frequ\u00eancia. Automodula\u00e7\u00e3o de fase. Polariza\u00e7\u00e3o das ondas eletromagn\u00e9ticas. Polariza\u00e7\u00e3o linear. Polariza\u00e7\u00e3o el\u00edptica. Polariza\u00e7\u00e3o circular. Obten\u00e7\u00e3o de luz linearmente polarizada. Equa\u00e7\u00f5es de Fresnel. Polariza\u00e7\u00e3o por reflex\u00e3o total interna. Matrizes de Jones. Atividade \u00f3ptica. Efeito Faraday. Isoladores \u00f3pticos. Efeito Pockels. Efeitos Kerr e Cotton-Mouton. Chaveamento eletro\u00f3ptico.Interfer\u00eancia. Princ\u00edpio da superposi\u00e7\u00e3o. Interfer\u00eancia por divis\u00e3o da frente de onda. Interfer\u00eancia por divis\u00e3o de amplitudes. Interfer\u00f4metro de Fabry-Perot. Analisador de espectro \u00f3ptico. Teoria de pel\u00edculas.Coer\u00eancia. Introdu\u00e7\u00e3o. Coer\u00eancia temporal. Resolu\u00e7\u00e3o espectral de um trem de ondas finito. Coer\u00eancia espacial. Medidas de di\u00e2metros de estrelas.Difra\u00e7\u00e3o. Princ\u00edpio de Huygens. F\u00f3rmula de Fresnel-Kirchhoff. Princ\u00edpio de Babinet. Difra\u00e7\u00e3o de Fraunhofer. Difra\u00e7\u00e3o por uma abertura circular. Rede de difra\u00e7\u00e3o. Padr\u00f5es de difra\u00e7\u00e3o de Fresnel. \u00d3ptica de Fourier.  Microscopia por contraste de fase.  Holografia. Intera\u00e7\u00e3o da radia\u00e7\u00e3o com a mat\u00e9ria.  Modelo do oscilador harm\u00f4nico.  Dispers\u00e3o crom\u00e1tica do \u00edndice de refra\u00e7\u00e3o. Absor\u00e7\u00e3o. Espalhamento Rayleigh. For\u00e7a da radia\u00e7\u00e3o em \u00e1tomo neutro.\u00d3ptica n\u00e3o linear. Susceptibilidade n\u00e3o linear, processos param\u00e9tricos e n\u00e3o param\u00e9tricos. Gera\u00e7\u00e3o de freq\u00fc\u00eancias. Casamento de fase.\" \"O que \u00e9 luz? Reflex\u00e3o. Refra\u00e7\u00e3o. Difra\u00e7\u00e3o. Polariza\u00e7\u00e3o. Forma\u00e7\u00e3o de imagens: Transformada de Fourier. Ondas eletromagn\u00e9ticas. Equa\u00e7\u00f5es de Maxwell. Propaga\u00e7\u00e3o da luz em diferentes meios: v\u00e1cuo, diel\u00e9trico, condutor. Transporte de energia. Condi\u00e7\u00f5es de contorno entre diferentes meios: v\u00e1cuo, diel\u00e9trico, condutor. Propaga\u00e7\u00e3o da luz entre diferentes meios: incid\u00eancia normal e obl\u00edqua na interface entre meios. Coeficientes de Fresnel. Aplica\u00e7\u00f5es da \u00d3ptica: holografia, laser, fibras \u00f3pticas, materiais eletrocr\u00f4micos, metamateriais.\" | Out-Null\nSet-ParagraphText $paragraphs \"Ray optics. Introduction. Propagation of light in homogeneous media. Propagation of light in non-homogeneous media. Generalized Snell's law. Fermat's principle. The equation of rays. The eikonal function. Analogy between classical mechanics and optics geometric The optical potential.Electromagnetic waves. One-dimensional harmonic waves. Flat and spherical waves. Gaussian waves. Gaussian beam propagation. Poynting vector. Intensity.The phase of the electromagnetic wave. Phase and group speeds. Dispersal. Doppler effect. Astronomical applications. Broadening of spectral lines. Relativistic optics. Electro-optical frequency modulation. Phase automodulation.Polarization of electromagnetic waves. Linear polarization. Elliptical Polarization. Circular polarization. Obtaining linearly polarized light. Fresnel equations. Polarization by total internal reflection. Jones matrices. Optical activity. Faraday effect. Optical isolators. Pockels Effect. Kerr and Cotton-Mouton effects. Electro-optical switching.Interference. Superposition principle. Interference by division of the wavefront. Amplitude division interference. Fabry-Perot interferometer. Optical spectrum analyzer. Film theory.Coherence. Introduction. Temporal coherence. Spectral resolution of a finite wave train. Spatial coherence. Star diameter measurements.Diffraction. Huygens Principle. Fresnel-Kirchhoff formula. Babinet's Principle. Fraunhofer Diffraction. Diffraction through a circular aperture. Diffraction grating. Fresnel diffraction patterns. Fourier optics. Phase contrast microscopy. Holography.Interaction of radiation with matter. Harmonic oscillator model. Chromatic dispersion of the refractive index. Absorption. Rayleigh scattering. Force of radiation on a neutral atom.Non-linear optics. Nonlinear susceptibility, parametric and nonparametric processes. Frequency generation. Phase marriage.\" \"What is light? Reflection. Refraction. Diffraction. Polarization. Image formation: Fourier transform. Electromagnetic waves. Maxwell's equations. Propagation of light in different media: vacuum, dielectric, conductor. Energy transport. Boundary conditions between different media: vacuum, dielectric, conductor. Light propagation between different media: normal and oblique incidence at the interface between media. Fresnel coefficients. Optics Applications: holography, laser, optical fibers, electrochromic materials, metamaterials.\" | Out-Null\n\n# 7) Append a new reference to the \"Bibliografia\" paragraph.\nfor ($i = 1; $i -le $paragraphs.Count; $i++) {\n    $p = $paragraphs.Item($i)\n    if ($p.Range.Text.TrimEnd(\"`r\", \"`a\") -eq \"HECHT, E.; ZAJAC, A. Optics; Reading, Addison-Wesley, 1974.ZILLIO, S. C. \u00d3ptica Moderna - Fundamentos e Aplica\u00e7\u00f5es, 2005.\") {\n        $p.Range.InsertAfter(\"J. R. Reitz, F. J. Milford, R. W. Christy, Fundamentos da Teoria Eletromagn\u00e9tica. Editora Campus. 1982.\")\n        break\n    }\n}\n\n"}
